$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7663920000000001
$ws.Range("H2").Value = 2.299176
$ws.Range("I2").Value = 0.4782815633346924
$ws.Range("J2").Value = 0.4782815633346925
$ws.Range("M2").Value = 43.19793166666667
$ws.Range("N2").Value = 129.593795
$ws.Range("O2").Value = 0.7412538312889448
$ws.Range("P2").Value = 0.7412538312889448
$ws.Range("Q2").Value = 33.10654924588
$ws.Range("R2").Value = 297.95894321292
$ws.Range("S2").Value = 0.3545280412567069
$ws.Range("T2").Value = 0.3545280412567069
$ws.Range("G3").Value = 0.7663920000000001
$ws.Range("H3").Value = 2.299176
$ws.Range("I3").Value = 0.4782815633346924
$ws.Range("J3").Value = 0.4782815633346925
$ws.Range("O3").Value = 0.1781209566020688
$ws.Range("P3").Value = 0.1781209566020688
$ws.Range("Q3").Value = 7.955399314720001
$ws.Range("R3").Value = 71.59859383248
$ws.Range("S3").Value = 0.08519196958630838
$ws.Range("T3").Value = 0.0851919695863084
$ws.Range("G4").Value = 0.7663920000000001
$ws.Range("H4").Value = 2.299176
$ws.Range("I4").Value = 0.4782815633346924
$ws.Range("J4").Value = 0.4782815633346925
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.422089
$ws.Range("N4").Value = 1.266267
$ws.Range("O4").Value = 0.007242825670663926
$ws.Range("P4").Value = 0.007242825670663927
$ws.Range("Q4").Value = 0.323485632888
$ws.Range("R4").Value = 2.911370695992
$ws.Range("S4").Value = 0.003464109984725784
$ws.Range("T4").Value = 0.003464109984725785
$ws.Range("G5").Value = 0.7663920000000001
$ws.Range("H5").Value = 2.299176
$ws.Range("I5").Value = 0.4782815633346924
$ws.Range("J5").Value = 0.4782815633346925
$ws.Range("M5").Value = 4.133026333333333
$ws.Range("N5").Value = 12.399079
$ws.Range("O5").Value = 0.07092056230936286
$ws.Range("P5").Value = 0.07092056230936288
$ws.Range("Q5").Value = 3.167518317656
$ws.Range("R5").Value = 28.507664858904
$ws.Range("S5").Value = 0.03391999741389753
$ws.Range("T5").Value = 0.03391999741389754
$ws.Range("G6").Value = 0.7663920000000001
$ws.Range("H6").Value = 2.299176
$ws.Range("I6").Value = 0.4782815633346924
$ws.Range("J6").Value = 0.4782815633346925
$ws.Range("M6").Value = 0.1434673333333333
$ws.Range("N6").Value = 0.430402
$ws.Range("O6").Value = 0.002461824128959449
$ws.Range("P6").Value = 0.002461824128959449
$ws.Range("Q6").Value = 0.109952216528
$ws.Range("R6").Value = 0.9895699487520001
$ws.Range("S6").Value = 0.001177445093053793
$ws.Range("T6").Value = 0.001177445093053793
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8359946666666667
$ws.Range("H7").Value = 2.507984
$ws.Range("I7").Value = 0.5217184366653075
$ws.Range("J7").Value = 0.5217184366653076
$ws.Range("M7").Value = 43.19793166666667
$ws.Range("N7").Value = 129.593795
$ws.Range("O7").Value = 0.7412538312889448
$ws.Range("P7").Value = 0.7412538312889448
$ws.Range("Q7").Value = 36.11324048436445
$ws.Range("R7").Value = 325.01916435928
$ws.Range("S7").Value = 0.3867257900322379
$ws.Range("T7").Value = 0.3867257900322379
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8359946666666667
$ws.Range("H8").Value = 2.507984
$ws.Range("I8").Value = 0.5217184366653075
$ws.Range("J8").Value = 0.5217184366653076
$ws.Range("O8").Value = 0.1781209566020688
$ws.Range("P8").Value = 0.1781209566020688
$ws.Range("Q8").Value = 8.677897731591111
$ws.Range("R8").Value = 78.10107958432
$ws.Range("S8").Value = 0.09292898701576044
$ws.Range("T8").Value = 0.09292898701576045
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8359946666666667
$ws.Range("H9").Value = 2.507984
$ws.Range("I9").Value = 0.5217184366653075
$ws.Range("J9").Value = 0.5217184366653076
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.422089
$ws.Range("N9").Value = 1.266267
$ws.Range("O9").Value = 0.007242825670663926
$ws.Range("P9").Value = 0.007242825670663927
$ws.Range("Q9").Value = 0.3528641528586667
$ws.Range("R9").Value = 3.175777375728
$ws.Range("S9").Value = 0.003778715685938141
$ws.Range("T9").Value = 0.003778715685938142
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8359946666666667
$ws.Range("H10").Value = 2.507984
$ws.Range("I10").Value = 0.5217184366653075
$ws.Range("J10").Value = 0.5217184366653076
$ws.Range("M10").Value = 4.133026333333333
$ws.Range("N10").Value = 12.399079
$ws.Range("O10").Value = 0.07092056230936286
$ws.Range("P10").Value = 0.07092056230936288
$ws.Range("Q10").Value = 3.455187971859555
$ws.Range("R10").Value = 31.096691746736
$ws.Range("S10").Value = 0.03700056489546532
$ws.Range("T10").Value = 0.03700056489546533
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8359946666666667
$ws.Range("H11").Value = 2.507984
$ws.Range("I11").Value = 0.5217184366653075
$ws.Range("J11").Value = 0.5217184366653076
$ws.Range("M11").Value = 0.1434673333333333
$ws.Range("N11").Value = 0.430402
$ws.Range("O11").Value = 0.002461824128959449
$ws.Range("P11").Value = 0.002461824128959449
$ws.Range("Q11").Value = 0.1199379255075556
$ws.Range("R11").Value = 1.079441329568
$ws.Range("S11").Value = 0.001284379035905656
$ws.Range("T11").Value = 0.001284379035905657
